# Update TC_ID Excel SCD0017 until SCD0025 and Update TC_ID Solution SCD0006 until SCD0025

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (tab) from SCD0274 to SCD0017
$ws.Name = "SCD0017"

# Update the TC_ID value in B2 from DGS-289 to SCD0017-004
$ws.Range("B2").Value = "SCD0017-004"

# Adjust column widths for B and C to fit new (longer) content
# (values chosen so the engine's pixel-quantized ColumnWidth lands as close as
# possible to the authored widths of 13.140625 / 39.5703125 characters)
$ws.Columns.Item(2).ColumnWidth = 12.3
$ws.Columns.Item(3).ColumnWidth = 38.65

# Update the view: zoom + selection
$ws.Application.ActiveWindow.Zoom = 93
$ws.Range("B3").Select()
